# Applies the "Updated cryptos list" edit described by the diff.
# Updates Price (D) and Volume(1h) (E) columns for the rows that changed,
# and swaps the Monero / PolygonEcosystemToken rows (38 <-> 39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.058.66"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "2.432.56"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.35"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.63"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "2.418.31"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.13"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.339"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.16"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "2.883.52"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "61.047.82"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "2.453.86"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.86"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.91"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.19"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.92"
$ws.Range("E27").Value = "  -3.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "574.33"
$ws.Range("E28").Value = "  -6.48%  "
$ws.Range("D29").Value = "2.555.17"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "0.0₃0913"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("E35").Value = "  -6.06%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.63"
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.31"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.369"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.30"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.12"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.73"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("E45").Value = "  -4.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.35"
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("D47").Value = "0.0₆0288"
$ws.Range("E47").Value = "  +26.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "141.50"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.53"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.593"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0507"
$ws.Range("E51").Value = "  -2.60%  "
